# Add benchmark results for "Ryzen 9 3950X" and "Athlon 200GE".
#
# The sheet is kept sorted (ascending) by column O ("26" header — actual
# seconds-per-iteration rating), matching the table's existing sortState.
# We insert the two new rows, populate them, then re-sort the whole data
# range so the new rows land in their correct position (same behaviour as
# Excel's Data > Sort that produced the original ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# --- 1) New row: Athlon 200GE -------------------------------------------
# Appended right after the current last data row; write the CPU-model
# name first so the new shared string "Athlon 200GE" is registered before
# the other two new strings (matches first-use ordering of the target
# workbook).
$athlonRow = $lastRow + 1
$ws.Range("A" + $athlonRow).Value = "AMD"
$ws.Range("B" + $athlonRow).Value = 1
$ws.Range("C" + $athlonRow).Value = "Athlon 200GE"
$ws.Range("D" + $athlonRow).Value = 35
$ws.Range("E" + $athlonRow).Value = 2
$ws.Range("F" + $athlonRow).Value = 2
$ws.Range("G" + $athlonRow).Value = 3.2
$ws.Range("H" + $athlonRow).Value = 3.2
$ws.Range("I" + $athlonRow).Value = "x86-64"
$ws.Range("K" + $athlonRow).Value = 8
$ws.Range("L" + $athlonRow).Value = 2
$ws.Range("M" + $athlonRow).Value = "DDR4"
$ws.Range("N" + $athlonRow).Value = 3000
$ws.Range("O" + $athlonRow).Value = 1.19
$ws.Range("P" + $athlonRow).Value = 2.48
$ws.Range("Q" + $athlonRow).Value = 5.13
$ws.Range("R" + $athlonRow).Value = 10.52

# Copy number formatting etc. from the row above so the new row matches
# the rest of the table (columns F/G/H carry a dedicated style).
$ws.Range("A" + ($athlonRow - 1) + ":S" + ($athlonRow - 1)).Copy()
$ws.Range("A" + $athlonRow + ":S" + $athlonRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# re-write the values (PasteSpecial(xlPasteFormats) only touched formats,
# but guard against any accidental clobber by re-asserting the data)
$ws.Range("A" + $athlonRow).Value = "AMD"
$ws.Range("B" + $athlonRow).Value = 1
$ws.Range("C" + $athlonRow).Value = "Athlon 200GE"
$ws.Range("D" + $athlonRow).Value = 35
$ws.Range("E" + $athlonRow).Value = 2
$ws.Range("F" + $athlonRow).Value = 2
$ws.Range("G" + $athlonRow).Value = 3.2
$ws.Range("H" + $athlonRow).Value = 3.2
$ws.Range("I" + $athlonRow).Value = "x86-64"
$ws.Range("K" + $athlonRow).Value = 8
$ws.Range("L" + $athlonRow).Value = 2
$ws.Range("M" + $athlonRow).Value = "DDR4"
$ws.Range("N" + $athlonRow).Value = 3000
$ws.Range("O" + $athlonRow).Value = 1.19
$ws.Range("P" + $athlonRow).Value = 2.48
$ws.Range("Q" + $athlonRow).Value = 5.13
$ws.Range("R" + $athlonRow).Value = 10.52
$ws.Range("J" + $athlonRow).Value = ""
$ws.Range("S" + $athlonRow).Value = ""

$lastRow = $athlonRow

# --- 2) New row: Ryzen 9 3950X ------------------------------------------
# Inserted right under row 2 (it is AMD's newest high core-count part,
# sitting immediately below the 5950X at the top of the table); pushes
# every following row down by one.
$ws.Rows.Item(3).Insert()
$newRow = 3

$ws.Range("A" + $newRow).Value = "AMD"
$ws.Range("B" + $newRow).Value = 1
$ws.Range("D" + $newRow).Value = 105
$ws.Range("E" + $newRow).Value = 16
$ws.Range("F" + $newRow).Value = 2
$ws.Range("G" + $newRow).Value = 3.5
$ws.Range("H" + $newRow).Value = 4.7
$ws.Range("I" + $newRow).Value = "x86-64"
$ws.Range("J" + $newRow).Value = "Ryzen Balanced"
$ws.Range("C" + $newRow).Value = "Ryzen 9 3950X"
$ws.Range("K" + $newRow).Value = 64
$ws.Range("L" + $newRow).Value = 2
$ws.Range("M" + $newRow).Value = "DDR4"
$ws.Range("N" + $newRow).Value = 3200
$ws.Range("O" + $newRow).Value = 0.22
$ws.Range("P" + $newRow).Value = 0.44
$ws.Range("Q" + $newRow).Value = 0.86
$ws.Range("R" + $newRow).Value = 1.76

$lastRow = $lastRow + 1

# --- 3) Re-sort the data range by column O (ascending), same key as the
#        table's existing sortState/sortCondition. -----------------------
$dataRange = $ws.Range("A2:S" + $lastRow)
$keyRange = $ws.Range("O2:O" + $lastRow)
$dataRange.Sort($keyRange)
